$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Helper: set a cell's value while forcing it to be stored as text, even if
# it looks like a number (e.g. "589.19"), matching the source data which is
# always plain text in this sheet. Resets formatting/style back afterwards
# so no visible formatting change is left behind.
function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# Each entry: row number, optional B (coin name), optional C (link),
# optional D (price), optional E (volume 1h)
$updates = @(
    @{ Row=2;  D="64.349.55"; E="  +0.12%  " },
    @{ Row=3;  D="3.497.00";  E="  +0.20%  " },
    @{ Row=4;  E="  +0.01%  " },
    @{ Row=5;  D="589.19";    E="  +0.43%  " },
    @{ Row=6;  D="134.27";    E="  +0.13%  " },
    @{ Row=8;  E="  +0.67%  " },
    @{ Row=9;  E="  +5.97%  " },
    @{ Row=10; E="  +0.55%  " },
    @{ Row=11; D="0.391";     E="  +3.73%  " },
    @{ Row=12; D="4.094.78";  E="  +0.18%  " },
    @{ Row=13; E="  +0.58%  " },
    @{ Row=14; E="  +0.53%  " },
    @{ Row=15; D="3.498.08";  E="  +0.07%  " },
    @{ Row=16; D="64.355.92"; E="  +0.05%  " },
    @{ Row=17; D="25.47";     E="  +1.40%  " },
    @{ Row=18; D="10.03";     E="  +0.43%  " },
    @{ Row=19; D="5.78";      E="  +0.70%  " },
    @{ Row=20; D="13.54";     E="  -0.49%  " },
    @{ Row=21; D="388.70";    E="  +0.39%  " },
    @{ Row=22; E="  +3.03%  " },
    @{ Row=23; D="3.637.18";  E="  +0.17%  " },
    @{ Row=24; D="74.29";     E="  -0.42%  " },
    @{ Row=25; E="  +0.11%  " },
    @{ Row=26; E="  -0.93%  " },
    @{ Row=27; E="  +2.76%  " },
    @{ Row=28; D="0.999";     E="  -0.06%  " },
    @{ Row=29; D="7.39" },
    @{ Row=30; D="2.26";      E="  +1.52%  " },
    @{ Row=31; D="1.49";      E="  -3.45%  " },
    @{ Row=32; D="8.17";      E="  -1.07%  " },
    @{ Row=33; E="  +5.61%  " },
    @{ Row=34; D="3.525.84";  E="  +0.43%  " },
    @{ Row=35; E="  +0.00%  " },
    @{ Row=36; D="23.37";     E="  -0.31%  " },
    @{ Row=37; D="5.34";      E="  +1.87%  " },
    @{ Row=38; E="  +1.58%  " },
    @{ Row=39; E="  +1.30%  " },
    @{ Row=40; D="165.62";    E="  +2.55%  " },
    @{ Row=41; D="0.0788";    E="  +1.01%  " },
    @{ Row=42; E="  +0.59%  " },
    @{ Row=43; E="  -0.01%  " },
    @{ Row=44; D="4.43";      E="  +0.87%  " },
    @{ Row=45; D="24.80";     E="  -2.31%  " },
    @{ Row=46; E="  +0.83%  " },
    @{ Row=47; D="1.65";      E="  -0.33%  " },
    @{ Row=48; B="Cosmos";     C="https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom";     D="6.83";   E="  +1.53%  " },
    @{ Row=49; B="SuiNetwork"; C="https://coinranking.com/coin/3xJluUMvp+suinetwork-sui";       D="0.923";  E="  +2.98%  " },
    @{ Row=50; D="2.405.17";  E="  -2.68%  " },
    @{ Row=51; E="  +0.01%  " }
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($u.ContainsKey("B")) { $ws.Cells.Item($r, 2).Value = $u.B }
    if ($u.ContainsKey("C")) { $ws.Cells.Item($r, 3).Value = $u.C }
    if ($u.ContainsKey("D")) { Set-TextValue $r 4 $u.D }
    if ($u.ContainsKey("E")) { $ws.Cells.Item($r, 5).Value = $u.E }
}
